$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 817-818, shifting the existing rows 817-883
# down to 819-885.
$ws.Range("A817:A818").EntireRow.Insert()

# Populate the first newly-inserted row (817) with its new record.
$ws.Range("A817").Value = 9
$ws.Range("B817").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C817").Value = "Metropolitana"
$ws.Range("D817").Value = 45265
$ws.Range("E817").Value = 13
$ws.Range("F817").Value = 100112031
$ws.Range("G817").Value = "Poroto verde"
$ws.Range("H817").Value = "Magnum"
$ws.Range("I817").Value = "Primera"
$ws.Range("J817").Value = 70
$ws.Range("K817").Value = 43000
$ws.Range("L817").Value = 45000
$ws.Range("M817").Value = 44000
$ws.Range("N817").Value = "$/malla 25 kilos"
$ws.Range("O817").Value = "Provincia de Limarí"
$ws.Range("P817").Value = 1760
$ws.Range("Q817").Value = 25
$ws.Range("R817").Value = "Hortaliza"

# Populate the second newly-inserted row (818) with its new record.
$ws.Range("A818").Value = 9
$ws.Range("B818").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C818").Value = "Metropolitana"
$ws.Range("D818").Value = 45265
$ws.Range("E818").Value = 13
$ws.Range("F818").Value = 100112031
$ws.Range("G818").Value = "Poroto verde"
$ws.Range("H818").Value = "Sin especificar"
$ws.Range("I818").Value = "Primera"
$ws.Range("J818").Value = 16
$ws.Range("K818").Value = 40000
$ws.Range("L818").Value = 41000
$ws.Range("M818").Value = 40500
$ws.Range("N818").Value = "$/malla 25 kilos"
$ws.Range("O818").Value = "Provincia de Huasco"
$ws.Range("P818").Value = 1620
$ws.Range("Q818").Value = 25
$ws.Range("R818").Value = "Hortaliza"
